# Generate Report for Archive
# Updates the localization status from "Ready for handoff" to "In Translation"
# on the Overview sheet (zh-cn/de-de status columns) and on each per-language
# status sheet, then re-fits the affected status columns to their new content.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: status is reported in columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$usedRows = $wsOverview.UsedRange.Rows.Count
for ($r = 2; $r -le $usedRows; $r++) {
    foreach ($col in @("E", "F")) {
        $cell = $wsOverview.Range($col + $r)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value2 = $newStatus
        }
    }
}
$wsOverview.Columns.Item("E").EntireColumn.AutoFit() | Out-Null
$wsOverview.Columns.Item("F").EntireColumn.AutoFit() | Out-Null

# --- Per-language sheets: status lives in column C ("Status") ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $usedRows = $ws.UsedRange.Rows.Count
    for ($r = 2; $r -le $usedRows; $r++) {
        $cell = $ws.Range("C" + $r)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value2 = $newStatus
        }
    }
    $ws.Columns.Item("C").EntireColumn.AutoFit() | Out-Null
}
